# 📊 Horarios actualizados Línea 141 - 721
# Updates the "last updated" timestamp, row-count summary, and appends the
# newly scraped arrival rows for each route sheet.

$wb = $excel.ActiveWorkbook

$newTimestamp = "02:24:16"

# ---------------------------------------------------------------------------
# Sheet 1: LP1912  -> gains 3 new rows (11-13), Total filas 5 -> 8
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2, 1).Value = "Última actualización: $newTimestamp"
$ws1.Cells.Item(3, 1).Value = "Total filas: 8"

$ws1.Cells.Item(11, 1).Value = $newTimestamp
$ws1.Cells.Item(11, 2).Value = "03:53"
$ws1.Cells.Item(11, 3).Value = "14_ABASTO"
$ws1.Cells.Item(11, 4).Value = 89
$ws1.Cells.Item(11, 5).Value = "LP1912"

$ws1.Cells.Item(12, 1).Value = $newTimestamp
$ws1.Cells.Item(12, 2).Value = "03:58"
$ws1.Cells.Item(12, 3).Value = "215_ALUAR"
$ws1.Cells.Item(12, 4).Value = 94
$ws1.Cells.Item(12, 5).Value = "LP1912"

$ws1.Cells.Item(13, 1).Value = $newTimestamp
$ws1.Cells.Item(13, 2).Value = "04:01"
$ws1.Cells.Item(13, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(13, 4).Value = 97
$ws1.Cells.Item(13, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215 -> gains 1 new row (9), Total filas 3 -> 4
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = "Última actualización: $newTimestamp"
$ws2.Cells.Item(3, 1).Value = "Total filas: 4"

$ws2.Cells.Item(9, 1).Value = $newTimestamp
$ws2.Cells.Item(9, 2).Value = "03:58"
$ws2.Cells.Item(9, 3).Value = "215_ALUAR"
$ws2.Cells.Item(9, 4).Value = 94
$ws2.Cells.Item(9, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173 -> only the "last updated" stamp changes
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = "Última actualización: $newTimestamp"
